$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that held the ad-hoc "12TRC1234 / TEST / MESSING" pretest
# scratch data was row 2; deleting it shifts every subsequent row up by
# one and drops the shared-string entries that were only referenced by
# that row.
$ws.Rows.Item(2).Delete()

# Re-select the header-adjacent row (mirrors Excel resetting the active
# cell to column A after a full-row delete).
$ws.Range("A2:XFD2").Select() | Out-Null
